$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 428.42856
$ws.Range("I2").Value = 433.16666
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 433.16666
$ws.Range("L2").Value = 400
$ws.Range("M2").Value = -320.16666
$ws.Range("N2").Value = -626
$ws.Range("H64").Value = 3297.5476
$ws.Range("I64").Value = 3071.2144
$ws.Range("J64").Value = 3410.7144
$ws.Range("K64").Value = 3071.2144
$ws.Range("L64").Value = 3410.7144
$ws.Range("M64").Value = -2823.2144
$ws.Range("N64").Value = -3906.7144
$ws.Range("H67").Value = 3297.5476
$ws.Range("I67").Value = 3071.2144
$ws.Range("J67").Value = 3410.7144
$ws.Range("K67").Value = 3071.2144
$ws.Range("L67").Value = 3410.7144
$ws.Range("M67").Value = -2213.2144
$ws.Range("N67").Value = -5126.7144
$ws.Range("H70").Value = 7985742.5
$ws.Range("I70").Value = 20960262
$ws.Range("J70").Value = 1423
$ws.Range("K70").Value = 62880786
$ws.Range("L70").Value = 4269
$ws.Range("M70").Value = -62880516
$ws.Range("N70").Value = -4809
$ws.Range("H73").Value = 7985742.5
$ws.Range("I73").Value = 20960262
$ws.Range("J73").Value = 1423
$ws.Range("K73").Value = 62880786
$ws.Range("L73").Value = 4269
$ws.Range("M73").Value = -62879850
$ws.Range("N73").Value = -6141
$ws.Range("H132").Value = 2206.9792
$ws.Range("I132").Value = 1708.4054
$ws.Range("J132").Value = 3884
$ws.Range("K132").Value = 5125.216200000001
$ws.Range("L132").Value = 11652
$ws.Range("M132").Value = -2595.216200000001
$ws.Range("N132").Value = -16712
$ws.Range("H135").Value = 1046.1724
$ws.Range("I135").Value = 873.5599999999999
$ws.Range("J135").Value = 2125
$ws.Range("K135").Value = 7862.039999999999
$ws.Range("L135").Value = 19125
$ws.Range("M135").Value = -5327.039999999999
$ws.Range("N135").Value = -24195
$ws.Range("H138").Value = 2556.6704
$ws.Range("I138").Value = 1609.1034
$ws.Range("J138").Value = 4388.6333
$ws.Range("K138").Value = 4827.3102
$ws.Range("L138").Value = 13165.8999
$ws.Range("M138").Value = 312.6898000000001
$ws.Range("N138").Value = -23445.8999
$ws.Range("H141").Value = 4930.3784
$ws.Range("I141").Value = 2267.2424
$ws.Range("J141").Value = 26901.25
$ws.Range("K141").Value = 6801.7272
$ws.Range("L141").Value = 80703.75
$ws.Range("M141").Value = -1621.7272
$ws.Range("N141").Value = -91063.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1354.7693
$ws.Range("I74").Value = 1379.1111
$ws.Range("J74").Value = 1300
$ws.Range("K74").Value = 1379.1111
$ws.Range("L74").Value = 1300
$ws.Range("M74").Value = -505.1111000000001
$ws.Range("N74").Value = -3048
$ws.Range("H77").Value = 1354.7693
$ws.Range("I77").Value = 1379.1111
$ws.Range("J77").Value = 1300
$ws.Range("K77").Value = 6895.5555
$ws.Range("L77").Value = 6500
$ws.Range("M77").Value = -2527.5555
$ws.Range("N77").Value = -15236

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1828.2858
$ws.Range("I99").Value = 964
$ws.Range("K99").Value = 964
$ws.Range("M99").Value = 534
$ws.Range("H132").Value = 68265.5
$ws.Range("J132").Value = 68265.5
$ws.Range("L132").Value = 68265.5
$ws.Range("N132").Value = -78385.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1735.9354
$ws.Range("I31").Value = 1303.0769
$ws.Range("J31").Value = 2469.913
$ws.Range("K31").Value = 1303.0769
$ws.Range("L31").Value = 2469.913
$ws.Range("M31").Value = -1008.0769
$ws.Range("N31").Value = -3059.913
$ws.Range("H34").Value = 1735.9354
$ws.Range("I34").Value = 1303.0769
$ws.Range("J34").Value = 2469.913
$ws.Range("K34").Value = 1303.0769
$ws.Range("L34").Value = 2469.913
$ws.Range("M34").Value = -1101.0769
$ws.Range("N34").Value = -2873.913
$ws.Range("H58").Value = 1685606.8
$ws.Range("I58").Value = 4631742
$ws.Range("J58").Value = 2100.7856
$ws.Range("K58").Value = 4631742
$ws.Range("L58").Value = 2100.7856
$ws.Range("M58").Value = -4631539
$ws.Range("N58").Value = -2506.7856
$ws.Range("H60").Value = 38900
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H132").Value = 376800.56
$ws.Range("I132").Value = 483773.78
$ws.Range("J132").Value = 2394.25
$ws.Range("K132").Value = 1451321.34
$ws.Range("L132").Value = 7182.75
$ws.Range("M132").Value = -1448791.34
$ws.Range("N132").Value = -12242.75
$ws.Range("H136").Value = 1685606.8
$ws.Range("I136").Value = 4631742
$ws.Range("J136").Value = 2100.7856
$ws.Range("K136").Value = 13895226
$ws.Range("L136").Value = 6302.3568
$ws.Range("M136").Value = -13892676
$ws.Range("N136").Value = -11402.3568

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 688.1739
$ws.Range("I122").Value = 500.33334
$ws.Range("J122").Value = 893.0909
$ws.Range("K122").Value = 4503.00006
$ws.Range("L122").Value = 8037.8181
$ws.Range("M122").Value = -2053.00006
$ws.Range("N122").Value = -12937.8181
$ws.Range("H132").Value = 2320.75
$ws.Range("I132").Value = 853.8
$ws.Range("J132").Value = 3368.5715
$ws.Range("K132").Value = 7684.2
$ws.Range("L132").Value = 30317.1435
$ws.Range("M132").Value = -5154.2
$ws.Range("N132").Value = -35377.1435

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 8000
$ws.Range("J47").Value = 8000
$ws.Range("L47").Value = 8000
$ws.Range("N47").Value = -9136
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("N65").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3072.0652
$ws.Range("I132").Value = 2746.9285
$ws.Range("J132").Value = 3577.8333
$ws.Range("K132").Value = 8240.7855
$ws.Range("L132").Value = 10733.4999
$ws.Range("M132").Value = -5710.7855
$ws.Range("N132").Value = -15793.4999
$ws.Range("H136").Value = 3372.44
$ws.Range("I136").Value = 2862.875
$ws.Range("J136").Value = 4278.3335
$ws.Range("K136").Value = 8588.625
$ws.Range("L136").Value = 12835.0005
$ws.Range("M136").Value = -6038.625
$ws.Range("N136").Value = -17935.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 948.67926
$ws.Range("I132").Value = 819.5238000000001
$ws.Range("J132").Value = 1441.8182
$ws.Range("K132").Value = 2458.5714
$ws.Range("L132").Value = 4325.4546
$ws.Range("M132").Value = 71.42859999999973
$ws.Range("N132").Value = -9385.454600000001

Write-Output "applied edits"